$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "Play Mega Moolah for Free - Exciting African Wildlife Slot Game" "Play Mega Moolah Free: Review and Gameplay"
Replace-Text "Four jackpot prizes with high payouts" "Four jackpots with high payout potential"
Replace-Text "Average volatility level for more consistent wins" "Average volatility level for more frequent wins"
Replace-Text "Bonus symbols and features increase payout opportunities" "Bonus symbols and free spin rounds"
Replace-Text "RTP is lower than other slots" "RTP of only 88.12%"
Replace-Text "Animal symbols may be too human-like for some players" "Some players may find the animal symbols unsettling"
Replace-Text "Experience the African wildlife-themed Mega Moolah slot game and play for free. With four jackpots, bonus symbols, and average volatility, it offers exciting payout opportunities." "Discover the features, pros, and cons of Mega Moolah and play this exciting slot game for free."
